# Update the two station pretty-names that changed on sheet "Blad1":
#   RS310 (row 41): Steenwijk -> Onna
#   RS406 (row 42): Eext      -> Gieten
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$ws.Range("B41").Value = "Onna"
$ws.Range("B42").Value = "Gieten"
